$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the "k" column (J) for the 10 data rows.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary labels (A) + aggregate formulas (B), with a bold,
# size 12, vertically centered style applied to the B column values.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page orientation (portrait, A4-ish "9" paper size) as set in the source file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on the newly added summary block, matching the saved view.
$ws.Range("A14:B17").Select() | Out-Null
